# Updates cryptos list with refreshed prices / volume(1h) percentages
# and restores the correct Coin/Link/Price ordering for a few rows
# that had shifted (BNB/USDC, Frax/HuobiToken, Quant/PaxDollar).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.243.71'
$ws.Range('E2').Value = '  -2.68%  '

$ws.Range('D3').Value = '1.777.26'
$ws.Range('E3').Value = '  -0.87%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.71'
$ws.Range('E5').Value = '  +1.58%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4219'
$ws.Range('E7').Value = '  +1.06%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3586'
$ws.Range('E8').Value = '  +1.09%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07088'
$ws.Range('E9').Value = '  +0.59%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8381'
$ws.Range('E10').Value = '  -0.19%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.37'
$ws.Range('E11').Value = '  +1.16%  '

$ws.Range('D12').Value = '1.740.30'
$ws.Range('E12').Value = '  -7.26%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.444'
$ws.Range('E13').Value = '  +1.80%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.214'
$ws.Range('E14').Value = '  -1.20%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.06924'
$ws.Range('E15').Value = '  +2.75%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.008'
$ws.Range('E16').Value = '  -0.09%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '78.87'
$ws.Range('E17').Value = '  -0.83%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008682'
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.13%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.90'
$ws.Range('E20').Value = '  -0.79%  '

$ws.Range('D21').Value = '26.246.43'
$ws.Range('E21').Value = '  -3.71%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.114'
$ws.Range('E22').Value = '  +1.24%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.97'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').Value = '1.938.40'
$ws.Range('E24').Value = '  -6.40%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '152.06'
$ws.Range('E25').Value = '  -0.52%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.803'
$ws.Range('E26').Value = '  -6.86%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.00'
$ws.Range('E27').Value = '  -0.48%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '5.060'
$ws.Range('E28').Value = '  +1.45%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '114.01'
$ws.Range('E29').Value = '  +0.96%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.844'
$ws.Range('E30').Value = '  +12.57%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08862'
$ws.Range('E31').Value = '  -0.51%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7247'
$ws.Range('E32').Value = '  +1.56%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.122'
$ws.Range('E33').Value = '  +4.95%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.306'
$ws.Range('E34').Value = '  +0.32%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.750'
$ws.Range('E35').Value = '  -3.48%  '

$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  -0.37%  '

$ws.Range('E37').Value = '  +3.41%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05098'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01881'
$ws.Range('E39').Value = '  -0.74%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.1608'
$ws.Range('E40').Value = '  -0.58%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4908'
$ws.Range('E41').Value = '  -0.60%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.588'
$ws.Range('E42').Value = '  +0.44%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '6.324'
$ws.Range('E43').Value = '  +5.04%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.060'
$ws.Range('E44').Value = '  +0.75%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '104.42'
$ws.Range('E45').Value = '  +0.22%  '

$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  -0.22%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.19'
$ws.Range('E47').Value = '  +0.35%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06168'
$ws.Range('E48').Value = '  -1.99%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.607'
$ws.Range('E49').Value = '  +1.20%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4446'
$ws.Range('E50').Value = '  -1.15%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.725'
$ws.Range('E51').Value = '  +3.96%  '
